$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top (shifts existing data down) and add a header label
$ws.Rows.Item(1).Insert()
$ws.Range("A1").Value = "HOUSE_LIST"

# Match the reported selection state
$ws.Range("C5").Select()
